# Refresh the "cryptos" price/volume table (Price = column D, Volume(1h) = column E)
# with the latest scraped values. Values that parse as plain numbers are forced
# back to text (leading "'" quote-prefix) so they keep matching the existing
# inlineStr/text formatting of this column (prices are stored as text, e.g.
# "26.037.84", "1.666.96", not numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.000.78'
$ws.Range("E2").Value = '  -2.13%  '
$ws.Range("D3").Value = '1.666.49'
$ws.Range("E3").Value = '  -1.50%  '
$ws.Range("D4").Value = "'" + '1.005'
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = "'" + '216.29'
$ws.Range("E5").Value = '  -1.86%  '
$ws.Range("D6").Value = "'" + '0.5095'
$ws.Range("E6").Value = '  -0.40%  '
$ws.Range("D7").Value = "'" + '1.005'
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("D8").Value = "'" + '0.2654'
$ws.Range("E8").Value = '  -0.80%  '
$ws.Range("E9").Value = '  +0.96%  '
$ws.Range("D10").Value = "'" + '21.81'
$ws.Range("E10").Value = '  -1.26%  '
$ws.Range("D11").Value = "'" + '0.07450'
$ws.Range("E11").Value = '  +1.08%  '
$ws.Range("D12").Value = '1.685.51'
$ws.Range("E12").Value = '  -0.57%  '
$ws.Range("D13").Value = "'" + '4.500'
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("D14").Value = "'" + '0.5819'
$ws.Range("E14").Value = '  +0.41%  '
$ws.Range("D15").Value = "'" + '0.000008529'
$ws.Range("E15").Value = '  -1.38%  '
$ws.Range("D16").Value = "'" + '64.07'
$ws.Range("E16").Value = '  -2.09%  '
$ws.Range("D17").Value = '26.144.17'
$ws.Range("E17").Value = '  -1.79%  '
$ws.Range("D18").Value = "'" + '4.925'
$ws.Range("E18").Value = '  -1.39%  '
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("E20").Value = '  -1.49%  '
$ws.Range("D21").Value = "'" + '190.74'
$ws.Range("E21").Value = '  +1.96%  '
$ws.Range("D22").Value = "'" + '6.187'
$ws.Range("E22").Value = '  -1.33%  '
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("D24").Value = "'" + '144.43'
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").Value = "'" + '7.607'
$ws.Range("E25").Value = '  +0.92%  '
$ws.Range("D26").Value = "'" + '0.1199'
$ws.Range("E26").Value = '  +1.78%  '
$ws.Range("E27").Value = '  -1.11%  '
$ws.Range("D28").Value = "'" + '0.06623'
$ws.Range("E28").Value = '  +14.18%  '
$ws.Range("D29").Value = "'" + '1.341'
$ws.Range("E29").Value = '  -0.59%  '
$ws.Range("E30").Value = '  -2.19%  '
$ws.Range("D31").Value = "'" + '3.553'
$ws.Range("E31").Value = '  +0.58%  '
$ws.Range("D32").Value = "'" + '3.512'
$ws.Range("E32").Value = '  -0.47%  '
$ws.Range("D33").Value = "'" + '1.658'
$ws.Range("E33").Value = '  -0.18%  '
$ws.Range("D34").Value = "'" + '1.016'
$ws.Range("E34").Value = '  -0.66%  '
$ws.Range("D35").Value = "'" + '0.6137'
$ws.Range("E35").Value = '  +2.24%  '
$ws.Range("D36").Value = "'" + '2.370'
$ws.Range("D37").Value = "'" + '2.689'
$ws.Range("E37").Value = '  +0.45%  '
$ws.Range("D38").Value = "'" + '6.341'
$ws.Range("E38").Value = '  +7.90%  '
$ws.Range("D39").Value = '1.095.36'
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("D40").Value = "'" + '0.01590'
$ws.Range("E40").Value = '  -2.08%  '
$ws.Range("D41").Value = "'" + '0.8694'
$ws.Range("E41").Value = '  +0.15%  '
$ws.Range("E42").Value = '  +0.26%  '
$ws.Range("D43").Value = "'" + '101.09'
$ws.Range("E43").Value = '  +1.35%  '
$ws.Range("D44").Value = '1.813.54'
$ws.Range("E44").Value = '  -1.98%  '
$ws.Range("D45").Value = "'" + '0.00000000109'
$ws.Range("E45").Value = '  -4.69%  '
$ws.Range("D46").Value = "'" + '56.33'
$ws.Range("E46").Value = '  -0.17%  '
$ws.Range("E47").Value = '  +0.28%  '
$ws.Range("D48").Value = "'" + '8.065'
$ws.Range("E48").Value = '  -0.21%  '
$ws.Range("D49").Value = "'" + '0.05232'
$ws.Range("E49").Value = '  -0.30%  '
$ws.Range("D50").Value = "'" + '0.4287'
$ws.Range("E50").Value = '  -0.83%  '
$ws.Range("D51").Value = "'" + '6.028'
$ws.Range("E51").Value = '  +3.07%  '
